$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 449-450; everything currently at/after row 449
# shifts down by two rows (so old row 449 -> 451, ..., old row 541 -> 543).
$ws.Rows("449:450").Insert()

# New record 1 (row 449)
$ws.Range("A449").Value = 11
$ws.Range("B449").Value = "Vega Monumental Concepción"
$ws.Range("C449").Value = "Bíobío"
$ws.Range("D449").Value = 44785
$ws.Range("E449").Value = 8
$ws.Range("F449").Value = 100112004
$ws.Range("G449").Value = "Cebolla"
$ws.Range("H449").Value = "Sin especificar"
$ws.Range("I449").Value = "1a (guarda)"
$ws.Range("J449").Value = 800
$ws.Range("K449").Value = 6000
$ws.Range("L449").Value = 6500
$ws.Range("M449").Value = 6250
$ws.Range("N449").Value = "`$/malla 18 kilos"
$ws.Range("O449").Value = "Región de O'Higgins"
$ws.Range("P449").Value = 347
$ws.Range("Q449").Value = 18
$ws.Range("R449").Value = "Hortaliza"

# New record 2 (row 450)
$ws.Range("A450").Value = 11
$ws.Range("B450").Value = "Vega Monumental Concepción"
$ws.Range("C450").Value = "Bíobío"
$ws.Range("D450").Value = 44785
$ws.Range("E450").Value = 8
$ws.Range("F450").Value = 100112004
$ws.Range("G450").Value = "Cebolla"
$ws.Range("H450").Value = "Sin especificar"
$ws.Range("I450").Value = "2a (guarda)"
$ws.Range("J450").Value = 400
$ws.Range("K450").Value = 5500
$ws.Range("L450").Value = 5500
$ws.Range("M450").Value = 5500
$ws.Range("N450").Value = "`$/malla 18 kilos"
$ws.Range("O450").Value = "Región de O'Higgins"
$ws.Range("P450").Value = 306
$ws.Range("Q450").Value = 18
$ws.Range("R450").Value = "Hortaliza"
